# "Add files via upload" — TODO list update.
# Three backlog items ("Permissions der Rolle ändern", "Nach Änderungen:
# Autorefresh", and "Rollen hinzufügen") move from in-progress (0%/0%/50%)
# to done (100%). Their "Soll" (D) description is replaced with a
# finished-tense sentence, and the old "Ist" (E) note describing the
# outstanding work is cleared out since the work is now complete.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 8: "Permissions der Rolle ändern" -> 100%, move the (reworded) note to D.
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "Die Permissions der Rolle können verändert werden"
$ws.Range("E8").Value = $null

# Row 10: "Nach Änderungen: Autorefresh" -> 100%, same treatment.
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "Nach jeder Änderung werden alle GUI-Tabellen aktualisiert"
$ws.Range("E10").Value = $null

# Row 7: "Rollen hinzufügen" -> 100%, describe what got added, drop old note.
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "Rolle hinzufügen können hinzugefügt werden"
$ws.Range("E7").Value = $null

# Selection moves from the old last-edited cell to the row just updated.
$ws.Range("D7").Select() | Out-Null
